# "severity sensation smuggle slope soak" added. this is last update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlTop vertical alignment constant (matches the other word/definition columns)
$xlTop = -4160

$words = @(
    @{ Row = 103; Height = 75;  Word = "severity";  Def = "severe problems, injuries, illnesses etc are very bad or very serious"; Ex1 = "he risk and severity of sunborn depend on he body's natural skin color."; Ex2 = "His injuries were quite severe." },
    @{ Row = 104; Height = 60;  Word = "sensation"; Def = "a feeling that you get from one of your five senses, especially the sense of touch"; Ex1 = "I experienced no sensation in my left foot."; Ex2 = "One sign of a heart attack is a tingling sensation in the left arm." },
    @{ Row = 105; Height = 45;  Word = "smuggle";   Def = "to take something or someone illegally from one country to another"; Ex1 = "if you try to smuggle drug you are stupid."; Ex2 = "The guns were smuggled across the border." },
    @{ Row = 106; Height = 75;  Word = "slope";     Def = "a surface of which one end or side is at a higher level than another; a rising or falling surface."; Ex1 = "the house builders slopped the roof..."; Ex2 = "the roof should have a slope sufficient for proper drainage" },
    @{ Row = 107; Height = 105; Word = "soak";      Def = "if you soak something, or if you let it soak, you keep it covered with a liquid for a period of time, especially in order to make it softer or easier to clean"; Ex1 = "Soak the clothes in cold water."; Ex2 = "soak the beans overnight in water" }
)

foreach ($w in $words) {
    $r = $w.Row
    $ws.Cells.Item($r, 1).Value = $w.Word
    $ws.Cells.Item($r, 1).VerticalAlignment = $xlTop
    $ws.Cells.Item($r, 2).Value = $w.Def
    $ws.Cells.Item($r, 3).Value = $w.Ex1
    $ws.Cells.Item($r, 4).Value = $w.Ex2
    $ws.Rows.Item($r).RowHeight = $w.Height
}

[void]$ws.Range("E107").Select()
